# Add a new row (row 4) with turma "1BADM" scheduling data, mirroring the
# layout of the existing rows (Turma, Segunda..Sexta counts).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A4:F4")

# Force the values to be stored as text (matching how the other rows' numeric
# looking entries, e.g. "1", "2", "3", are stored as text rather than numbers),
# then drop back to the default "Normal" style so no extra formatting/style is
# left behind on the new cells.
$newRow.NumberFormat = "@"

$ws.Range("A4").Value = "1BADM"
$ws.Range("B4").Value = "2"
$ws.Range("C4").Value = "1"
$ws.Range("D4").Value = "4"
$ws.Range("E4").Value = "2"
$ws.Range("F4").Value = "1"

$newRow.Style = "Normal"
